# Update SUMMARY_EXPORT_DATA_DOWNLOAD.xlsx data:
#  - "Short Term" sheet: revise rows 110-120 (cols B-G) and append new rows 121-124
#  - "Medium Term" sheet: revise rows 101-106 (cols B-D) and append new rows 107-110

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Short Term"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Short Term")

# Existing rows whose B:G values were revised.
$shortTermUpdates = @(
    @{ Row = 110; B = 7.93;  C = 15.16;  D = -12.49; E = 19.24;  F = 22.01;  G = -15.95 },
    @{ Row = 111; B = -4.26; C = -2.36;  D = 6.61;   E = -8.85;  F = -4.65;  G = -8.07 },
    @{ Row = 112; B = 2.74;  C = 2.66;   D = 4.91;   E = -2.58;  F = 2.31;   G = -0.63 },
    @{ Row = 113; B = 0.9;   C = 5.32;   D = -3.22;  E = 22.55;  F = 27.97;  G = -0.34 },
    @{ Row = 114; B = -2.26; C = -9.95;  D = -1.83;  E = 4.53;   F = 8.66;   G = -10.37 },
    @{ Row = 115; B = -2.94; C = -1.16;  D = 6.98;   E = -11.68; F = -8.94;  G = -6.89 },
    @{ Row = 116; B = 7.47;  C = 4.84;   D = 8.3;    E = 11.62;  F = 16.15;  G = 11.04 },
    @{ Row = 117; B = -0.02; C = 0.79;   D = -2.35;  E = 6.51;   F = 7.78;   G = 7.06 },
    @{ Row = 118; B = 17.82; C = 22.26;  D = -8.91;  E = 35.86;  F = 40;     G = -5.29 },
    @{ Row = 119; B = -8.31; C = -5.71;  D = -2.88;  E = 21.73;  F = 21.42;  G = -2.37 },
    @{ Row = 120; B = 4.69;  C = -2.46;  D = 7.56;   E = 26.86;  F = 31.34;  G = -16.15 }
)

foreach ($u in $shortTermUpdates) {
    $r = $u.Row
    $ws1.Cells.Item($r, 2).Value = $u.B
    $ws1.Cells.Item($r, 3).Value = $u.C
    $ws1.Cells.Item($r, 4).Value = $u.D
    $ws1.Cells.Item($r, 5).Value = $u.E
    $ws1.Cells.Item($r, 6).Value = $u.F
    $ws1.Cells.Item($r, 7).Value = $u.G
}

# New rows appended after the old last row (120).
$shortTermNewRows = @(
    @{ Row = 121; A = 45627; B = -8.36; C = -10.04; D = -3.07;  E = 15.27;  F = 17.27; G = -1.23 },
    @{ Row = 122; A = 45658; B = 33.67; C = 37.01;  D = -0.05;  E = 39.65;  F = 36.83; G = 10.82 },
    @{ Row = 123; A = 45689; B = 9.69;  C = 4.76;   D = 20.77;  E = 66.87;  F = 58.24; G = 22.15 },
    @{ Row = 124; A = 45717; B = 21.82; C = 34.64;  D = -11.9;  E = 101.78; F = 94.31; G = 5.43 }
)

foreach ($n in $shortTermNewRows) {
    $r = $n.Row
    # Copy the date cell's format (style index 1, numFmtId 14) from the row above
    # so the new date cell matches the existing column-A date styling.
    $ws1.Cells.Item($r - 1, 1).Copy()
    $ws1.Cells.Item($r, 1).PasteSpecial(-4122)

    $ws1.Cells.Item($r, 1).Value = $n.A
    $ws1.Cells.Item($r, 2).Value = $n.B
    $ws1.Cells.Item($r, 3).Value = $n.C
    $ws1.Cells.Item($r, 4).Value = $n.D
    $ws1.Cells.Item($r, 5).Value = $n.E
    $ws1.Cells.Item($r, 6).Value = $n.F
    $ws1.Cells.Item($r, 7).Value = $n.G
}

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Sheet 2: "Medium Term"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Medium Term")

# Existing rows whose B:D values were revised.
$mediumTermUpdates = @(
    @{ Row = 101; B = 8.05;  C = 7.03;  D = -0.82 },
    @{ Row = 102; B = 4.58;  C = 6.14;  D = 0.06 },
    @{ Row = 103; B = 4.26;  C = 8.2;   D = 2.3 },
    @{ Row = 104; B = 21.26; C = 14.48; D = 7.34 },
    @{ Row = 105; B = 22.93; C = 13.65; D = 9.25 },
    @{ Row = 106; B = 30.6;  C = 17.34; D = 12.39 }
)

foreach ($u in $mediumTermUpdates) {
    $r = $u.Row
    $ws2.Cells.Item($r, 2).Value = $u.B
    $ws2.Cells.Item($r, 3).Value = $u.C
    $ws2.Cells.Item($r, 4).Value = $u.D
}

# New rows appended after the old last row (106).
$mediumTermNewRows = @(
    @{ Row = 107; A = 45627; B = 23.58; C = 22.42; D = 14.48 },
    @{ Row = 108; A = 45658; B = 29.5;  C = 26.19; D = 16.02 },
    @{ Row = 109; A = 45689; B = 38.31; C = 34.42; D = 21.14 },
    @{ Row = 110; A = 45717; B = 63.55; C = 44.88; D = 29.96 }
)

foreach ($n in $mediumTermNewRows) {
    $r = $n.Row
    $ws2.Cells.Item($r - 1, 1).Copy()
    $ws2.Cells.Item($r, 1).PasteSpecial(-4122)

    $ws2.Cells.Item($r, 1).Value = $n.A
    $ws2.Cells.Item($r, 2).Value = $n.B
    $ws2.Cells.Item($r, 3).Value = $n.C
    $ws2.Cells.Item($r, 4).Value = $n.D
}

$excel.CutCopyMode = $false

Write-Output "SUMMARY_EXPORT_DATA_DOWNLOAD updates applied"
